$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1008.7692
$ws.Range("I32").Value = 1024.5
$ws.Range("J32").Value = 1005.9091
$ws.Range("K32").Value = 1024.5
$ws.Range("L32").Value = 1005.9091
$ws.Range("M32").Value = -698.5
$ws.Range("N32").Value = -1657.9091
$ws.Range("H40").Value = 6520
$ws.Range("J40").Value = 10099.8
$ws.Range("L40").Value = 10099.8
$ws.Range("N40").Value = -10449.8
$ws.Range("H64").Value = 7166.6665
$ws.Range("J64").Value = 7166.6665
$ws.Range("L64").Value = 7166.6665
$ws.Range("N64").Value = -7662.6665
$ws.Range("H67").Value = 7166.6665
$ws.Range("J67").Value = 7166.6665
$ws.Range("L67").Value = 7166.6665
$ws.Range("N67").Value = -8882.666499999999
$ws.Range("H74").Value = 11887.588
$ws.Range("I74").Value = 10577.857
$ws.Range("K74").Value = 10577.857
$ws.Range("M74").Value = -9641.857
$ws.Range("H77").Value = 11887.588
$ws.Range("I77").Value = 10577.857
$ws.Range("K77").Value = 52889.285
$ws.Range("M77").Value = -48209.285
$ws.Range("H132").Value = 10639.8
$ws.Range("I132").Value = 1897.3414
$ws.Range("J132").Value = 100250
$ws.Range("K132").Value = 5692.0242
$ws.Range("L132").Value = 300750
$ws.Range("M132").Value = -3162.0242
$ws.Range("N132").Value = -305810
$ws.Range("H137").Value = 6587.0435
$ws.Range("I137").Value = 6695.4443
$ws.Range("J137").Value = 6196.8
$ws.Range("K137").Value = 20086.3329
$ws.Range("L137").Value = 18590.4
$ws.Range("M137").Value = -17536.3329
$ws.Range("N137").Value = -23690.4
$ws.Range("H138").Value = 5773.7573
$ws.Range("J138").Value = 6156.0835
$ws.Range("L138").Value = 18468.2505
$ws.Range("N138").Value = -28748.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3426.9756
$ws.Range("I32").Value = 3187.4875
$ws.Range("J32").Value = 13006.5
$ws.Range("K32").Value = 3187.4875
$ws.Range("L32").Value = 13006.5
$ws.Range("M32").Value = -2900.4875
$ws.Range("N32").Value = -13580.5
$ws.Range("H45").Value = 2999
$ws.Range("I45").Value = 1998.4
$ws.Range("K45").Value = 1998.4
$ws.Range("M45").Value = -1621.4
$ws.Range("H61").Value = 3391.0356
$ws.Range("I61").Value = 3164.6667
$ws.Range("K61").Value = 3164.6667
$ws.Range("M61").Value = -2952.6667
$ws.Range("H74").Value = 1950.091
$ws.Range("I74").Value = 1791.6
$ws.Range("K74").Value = 1791.6
$ws.Range("M74").Value = -917.5999999999999
$ws.Range("H77").Value = 1950.091
$ws.Range("I77").Value = 1791.6
$ws.Range("K77").Value = 8958
$ws.Range("M77").Value = -4590
$ws.Range("H122").Value = 5511.2583
$ws.Range("I122").Value = 5429.615
$ws.Range("K122").Value = 16288.845
$ws.Range("M122").Value = -13838.845
$ws.Range("H123").Value = 84992
$ws.Range("J123").Value = 84992
$ws.Range("L123").Value = 84992
$ws.Range("N123").Value = -94792
$ws.Range("H135").Value = 80606.5
$ws.Range("J135").Value = 80606.5
$ws.Range("L135").Value = 80606.5
$ws.Range("N135").Value = -90746.5
$ws.Range("H136").Value = 3391.0356
$ws.Range("I136").Value = 3164.6667
$ws.Range("K136").Value = 9494.000100000001
$ws.Range("M136").Value = -6944.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 806.25
$ws.Range("I12").Value = 806.25
$ws.Range("K12").Value = 806.25
$ws.Range("M12").Value = -638.25
$ws.Range("H86").Value = 813002.9
$ws.Range("I86").Value = 1310972.8
$ws.Range("J86").Value = 3802
$ws.Range("K86").Value = 1310972.8
$ws.Range("L86").Value = 3802
$ws.Range("M86").Value = -1309849.8
$ws.Range("N86").Value = -6048
$ws.Range("H89").Value = 813002.9
$ws.Range("I89").Value = 1310972.8
$ws.Range("J89").Value = 3802
$ws.Range("K89").Value = 6554864
$ws.Range("L89").Value = 19010
$ws.Range("M89").Value = -6549248
$ws.Range("N89").Value = -30242
$ws.Range("H105").Value = 3200.2856
$ws.Range("I105").Value = 3200.2856
$ws.Range("K105").Value = 3200.2856
$ws.Range("M105").Value = -1453.2856
$ws.Range("H107").Value = 1252825
$ws.Range("I107").Value = 1797.8
$ws.Range("J107").Value = 3337870.2
$ws.Range("K107").Value = 1797.8
$ws.Range("L107").Value = 3337870.2
$ws.Range("M107").Value = 122.2
$ws.Range("N107").Value = -3341710.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 521.4286
$ws.Range("I22").Value = 441.66666
$ws.Range("K22").Value = 441.66666
$ws.Range("M22").Value = -91.66665999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 51069.383
$ws.Range("J107").Value = 88528.664
$ws.Range("L107").Value = 265585.992
$ws.Range("N107").Value = -269425.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1895.125
$ws.Range("I102").Value = 1060.1666
$ws.Range("K102").Value = 1060.1666
$ws.Range("M102").Value = 561.8334
$ws.Range("H132").Value = 42665.57
$ws.Range("I132").Value = 7110.25
$ws.Range("K132").Value = 21330.75
$ws.Range("M132").Value = -18800.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1212.375
$ws.Range("I22").Value = 1400
$ws.Range("J22").Value = 1149.8334
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 1149.8334
$ws.Range("M22").Value = -1105
$ws.Range("N22").Value = -1739.8334
$ws.Range("H27").Value = 1212.375
$ws.Range("I27").Value = 1400
$ws.Range("J27").Value = 1149.8334
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 1149.8334
$ws.Range("M27").Value = -1293
$ws.Range("N27").Value = -1363.8334
$ws.Range("H30").Value = 14012.5
$ws.Range("I30").Value = 12010.667
$ws.Range("J30").Value = 20018
$ws.Range("K30").Value = 12010.667
$ws.Range("L30").Value = 20018
$ws.Range("M30").Value = -11902.667
$ws.Range("N30").Value = -20234
$ws.Range("H46").Value = 2579.0833
$ws.Range("I46").Value = 2524.7
$ws.Range("K46").Value = 2524.7
$ws.Range("M46").Value = -2336.7
$ws.Range("H55").Value = 2432
$ws.Range("I55").Value = 220
$ws.Range("K55").Value = 220
$ws.Range("M55").Value = -47
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H82").Value = 1515.091
$ws.Range("I82").Value = 1911.8334
$ws.Range("K82").Value = 1911.8334
$ws.Range("M82").Value = -1550.8334
$ws.Range("H85").Value = 1515.091
$ws.Range("I85").Value = 1911.8334
$ws.Range("K85").Value = 1911.8334
$ws.Range("M85").Value = -663.8334
$ws.Range("H132").Value = 7325.737
$ws.Range("I132").Value = 6246.0386
$ws.Range("K132").Value = 18738.1158
$ws.Range("M132").Value = -16208.1158
$ws.Range("H138").Value = 58429
$ws.Range("J138").Value = 58429
$ws.Range("L138").Value = 58429
$ws.Range("N138").Value = -68709

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28932.426
$ws.Range("I132").Value = 1621.4546
$ws.Range("K132").Value = 4864.3638
$ws.Range("M132").Value = -2334.3638
